$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = '패션명품 쓱세일 COMING SOON'
$ws.Range("B3").Value = 'SSG머니 100% 페이백!'
$ws.Range("C3").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000004758'
$ws.Range("D3").Value = '(3/2~5) 패션명품 쓱세일 COMING SOON'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '2023-03-02'
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = '2023-03-05'
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").Value = '[''이벤트/쿠폰 > (3/2~5) 패션명품 쓱세일 COMING SOON'', ''스마일클럽'']'

# Row 4
$ws.Range("A4").Value = '첫구매는 반값다딜 <풀무원 편>'
$ws.Range("B4").Value = '풀무원 상품 4개 골라 최대 2만원 할인 + 무료배송'
$ws.Range("C4").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000004661&domainSiteNo=7018'
$ws.Range("D4").Value = '반값다딜 - 풀무원 편 (3/2~12)'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '2023-03-02'
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = '2023-03-12'
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").Value = '[''이벤트/쿠폰 > 반값다딜 - 풀무원 편 (3/2~12)'', ''스마일클럽'', ''STEP1 최근 1년간 쓱배송이 처음이라면, 웰컴 쿠폰 받기'', ''Pulmuone 상품쿠폰 50% 4장, 개별 쿠폰 당 최대 5천원 할인'', ''Pulmuone 상품쿠폰 무료배송 1장, 쓱/새벽배송 2만원 이상 구매 시'', ''      최근 1년 간 이마트몰, 트레이더스 쓱배송/점보택배 및 새벽배송 구매 이력이 없는 고객'', ''       쿠폰 발급 기간'', ''       쿠폰 사용 기간'', ''       상품할인쿠폰, 무료배송쿠폰: 2023년 3월 2일(목)~ 2023년 3월 12일 (일)'', ''       쿠폰 사용 조건'', ''-무료배송 쿠폰 : 이마트 쓱배송/새벽배송 상품 2만원 이상 구매시 무료배송'', ''       쿠폰 발급 대상'', ''(50%, 무료배송) [풀무원] 얇은피 꽉찬속 김치만두 400gx2 할인전 : 9,480원 할인후 : 4,740원'', ''(50%, 무료배송) [풀무원] 국산콩 투컵 두부 600g 할인전 : 5,480원 할인후 : 2,740원'', ''(50%, 무료배송) [풀무원] 바릴라 캠핑 기획팩 2,700G 할인전 : 29,800원 할인후 : 14,900원'', ''(50%, 무료배송) [풀무원] 새콤달콤 유부초밥 330g 할인전 : 5,680원 할인후 : 2,840원'', ''(50%, 무료배송) [풀무원] 행복한 동물복지유정란 15구 780g : 8,980원 할인후 : 4,490원'', ''할인 상품 한 눈에 보러가기'']'

# Row 5
$ws.Range("A5").Value = '어서와 봄, 할인상품 득템해 봄'
$ws.Range("B5").Value = '창고 대개방! 반짝할인 놓치지 말고 득템하세요'
$ws.Range("C5").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000004717'
$ws.Range("D5").Value = '(3/2~31) 창고대개방! 어서 와 봄, 할인상품 득템해 봄'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '2023-03-02'
$ws.Range("E5").Style = "Normal"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = '2023-03-31'
$ws.Range("F5").Style = "Normal"
$ws.Range("G5").Value = '[''이벤트/쿠폰 > (3/2~31) 창고대개방! 어서 와 봄, 할인상품 득템해 봄'', ''스마일클럽'', ''(3/2~31) 창고대개방! 어서 와 봄, 할인상품 득템해 봄'']'

# Row 6
$ws.Range("A6").Value = '최대 7.2만원 혜택'
$ws.Range("B6").Value = '+ 쓸 때마다 최대 12% 적립'
$ws.Range("C6").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000000858&siteNo=6005&recruitmentPath=L6007001&eventCode=HPG02'
$ws.Range("D6").Value = 'SSG.COM카드 Edition 2 이벤트 안내 페이지'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '2022-07-08'
$ws.Range("E6").Style = "Normal"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = '2023-06-28'
$ws.Range("F6").Style = "Normal"
$ws.Range("G6").Value = '[''이벤트/쿠폰 > SSG.COM카드 Edition 2 이벤트 안내 페이지'', ''스마일클럽'', ''SSG.COM카드 Edition 2 이벤트 안내 페이지'', ''SSG.COM카드 Edition2는 SSG.COM에서 최대 7.2만원 혜택 + 쓸 때마다 최대 12% 적립'', ''이달의 혜택 01. SSG.COM카드 Edition2 첫결제 시 5.2만원 쿠폰 할인 바로보기'', ''혜택 01. SSG머니 최대 12% 적립 바로보기'', ''혜택 02. 매월 스마일클럽 가입비 3,900원 지원 바로보기'', ''이달의 혜택. SSG.COM에서 SSG.COM카드 Edition2 첫 결제 시 5.2만원 할인'', ''      직전 6개월간 (2022년 9월 1일 부터 2023년 2월 28일) SSG.COM카드 및 SSG.COM카드 Edition2로 결제 이력이 없고 & SSGPAY에 등록된 현대카드가 없는 회원에 한함'', ''첫 결제 쿠폰 이용 방법'', ''② SSGPAY에 등록된 SSG.COM카드 Edition2로'', ''③ SSG.COM에서 기간 내 5만 3천원 이상 첫 결제 시 5.2만원 쿠폰 할인'', ''      본 혜택은 SSG.COM카드 및 SSG.COM카드 Edition2 로 SSG.COM(이마트몰, 신세계몰, 신세계백화점몰 등)에서 직전 6개월간 (2022.08.01~2023.01.31) SSG.COM카드 및 SSG.COM카드 Edition2로 결제 이력이 없고 & SSGPAY에 등록된 현대카드가 없는 회원에 한해 제공됨'', "      본 혜택은 결제 시 [결제방법 > SSGPAY카드] 내 SSG.COM카드 Edition2 선택 시, ''카드할인 최적가'' 추천에 의해 할인 금액이 자동 적용됨. 단, 5만 3천원 이상 첫 결제 시 적용)", ''다운받은 쿠폰은 SSGPAY에 등록된 SSG.COM카드 Edition2로 SSG.COM에서 바로 결제 시 사용 가능합니다.'', ''쿠폰 사용하러 가기'', ''SSGPAY에 등록된 SSG.COM카드 Edition2로 SSG.COM에서 바로 결제 시 사용 가능합니다.'', ''혜택 01. 장 볼 때마다 SSG머니 최대 12% 적립'', ''SSG.COM에서 최대 12%'', ''       스마일클럽으로 5% 적립'', ''       쓱·새벽·트레이더스 구매 시 (구매 전 스마일클럽 적립 아이콘을 꼭 확인해주세요)'', ''       SSG.COM 카드 Edition2로 7% 적립'', ''어디서나 한도없이 0.5%'', ''SSG.COM카드 Edition2로 어디서나 한도없이 0.5% 적립 (SSG.COM 외 모든 가맹점)'', ''SSG Money 최대 12% 적립'', ''스마일클럽 5% 적립 + SSG.COM카드 Edition2 최대 7% 적립'', ''      스마일클럽 5% 적립은 쓱·새벽배송·트레이더스 이용 시에 한함'', ''      SSG.COM카드 Edition2 최대 7% 적립'', ''      SSG.COM에서 결제 시 7%(1만 쓱머니 한도), 그 외 가맹점 0.5% 적립(적립한도 제한 없음)'', ''      무이자 할부 및 현대카드에서 제공하는 다른 할인 서비스 이용 시 적립 제외'', ''스마일클럽 가입비 3,900원 매월 100% 지원'', ''      1. 스마일클럽 자동 가입에'', ''      2. SSGPAY 내 카드 자동 등록'', ''      3. 스마일클럽 정기결제수단 자동 등록 및 월 이용료 3,900원 지원까지! (단, 해당 카드를 월 정기결제 수단에 등록한 경우에 한함)'', ''SSG.COM카드 Edition2를 스마일클럽 월 정기결제 수단에 등록 및 전월 이용금액 30만원 이상 시 혜택 제공'', ''TIP. 스마일클럽 가입 시 SSG.COM 혜택'', ''(SSG 가입 시) 15,000원만 담아도 쓱 · 새벽배송 무료배송'', ''쓱 · 새벽배송 · 트레이더스 SSG머니 최대 5% 적립'', ''장바구니 최대 10% 할인쿠폰'', ''매월 4장씩 최대 12% 할인쿠폰'', ''스마일배송 1만 5천원 이상 무제한 무료배송'', ''스마일배송 상품 스마일캐시 1% 적립'', ''스마일클럽 단독 혜택 스타벅스 상품 전용 딜'', ''스마일클럽 가입비 매월 3,900원 지원'', ''월 1회, 매달 스마일클럽 정기결제일에 혜택 제공'', ''SSG.COM카드 Edition2는 최초 발급 시, 스마일클럽 월 정기결제 수단에 자동 등록 됨'', ''전월 이용금액 30만원 미만 시, SSG.COM카드 Edition2로 스마일클럽 정기 결제 금액이 자동 결제됨'', ''스마일클럽 무료 이용 기간이라면 정기결제 금액 지원 대신 SSG머니 3,900원 제공'', ''      스타벅스 자동 충전, 생활요금(통신요금, 아파트관리비 등) 정기결제 신청 및 이체 시 최대 1만원 청구 할인'', ''      2. 스타벅스 자동 충전 또는 생활요금 정기결제 신청(각 항목당 할인한도 5천원, 최대 1만원 할인)'', ''      이용금액이 혜택금액보다 적을 경우, 이용 금액만큼 할인 적용'', ''      정기결제 신청 후 카드 결제일에 따라 매출 발생 다음 달 또는 다다음 달 청구 할인 혜택 적용'', ''      단, 청구 할인 제공 일정은 당사 또는 신청인 사정에 의해 상이할 수 있음'', ''      3. 쏘카 1만원 할인쿠폰'', ''      쿠폰은 등록일 포함 30일간 이용 가능'', ''실물 SSG.COM카드 Edition2 수령 전 SSGPAY로 결제 시 건당 100만원 이하 결제 가능 *단, 본인 확인(신분증 확인 및 1원 인증) 완료한 경우에 한하며, 건당 100만원 초과 시 실물카드 수령 후 결제 가능'', ''SSG머니 최대 적립 12%에서 5%는 SSG.COM에서 제공하는 멤버십 서비스로 SSG.COM 사정에 따라 변경 가능함'', ''카드 이용대금 연체 시 약정금리 + 연체가산금리 3%의 연체이자율이 적용됩니다. (회원별, 이용 상품별 차등적용 / 법정 최고금리 20% 이내) 단, 연체 발생시점에 약정금리가 없는 경우 아래와 같이 적용'', ''일시불 : 거래 발생시점 기준 최소 기간 (2개월)의 유이자 할부 약정금리 + 연체가산금리 3%'', ''무이자할부 : 거래발생시점 기준 동일한 할부 계약 기간의 유이자할부 약정금리 + 연체가산금리 3%'']'

# Row 7
$ws.Range("A7").Value = '더 강력해진 SSG.COM 삼성카드'
$ws.Range("B7").Value = 'SSG MONEY 최대 15% 적립 + 스마일클럽 월이용료 할인'
$ws.Range("C7").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000002385&recruitmentPath=SSG'
$ws.Range("D7").Value = 'SSG.COM삼성카드 리뉴얼 이벤트 안내 페이지'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '2022-10-26'
$ws.Range("E7").Style = "Normal"
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = '2025-10-25'
$ws.Range("F7").Style = "Normal"
$ws.Range("G7").Value = '[''이벤트/쿠폰 > SSG.COM삼성카드 리뉴얼 이벤트 안내 페이지'', ''스마일클럽'', ''SSG.COM삼성카드 리뉴얼 이벤트 안내 페이지'', ''SSG머니 최대 15% 적립'', ''스마일클럽 월 이용료 3,900원 결제일 할인'', ''최대 15% SSG머니 적립 중 5%는 SSG.COM에서 제공하는 서비스로 자세한 내용은 SSG.COM 확인바람.'', ''01. SSG MONEY 최대 15% 적립!'', ''       이용실적 관계없이 적립 : 1 % + 이용실적 충족 시 적립 4% = 적립 가능한 최대 SSG머니 5%'', ''       카드혜택 + 스마일클럽 혜택'', ''       이용실적 관계없이 적립 1% + 이용실적 충족 시 적립 9% + 쓱배송/쓱배송 traders/새벽배송 상품 구매 시 5% = 적립 가능한 최대 SSG머니 15%'', ''최대 15% SSG MONEY 적립'', ''SSG.COM 삼성카드 최대 10% 적립+스마일클럽 5% 적립'', ''SSG.COM 삼성카드 최대 10% 적립(1% 적립+9% 추가 적립)'', ''1% 적립(전월 실적 조건 및 적립한도 없음)'', ''9% 추가 적립'', ''삼성카드 할인이 적용된 일시불 및 할부 이용금액은 제외됩니다.'', ''스마일클럽 5% 적립'', ''쓱배송/쓱배송 traders/새벽배송 상품 결제건에 한해 혜택을 받을 수 있습니다.'', ''02. 스마일클럽 월이용료 매월 3,900원 결제일 할인'', ''15,000원만 담아도 쓱 새벽배송 무료배송(SSG가입시)'', ''1쓱 새벽배송 트레이더 - SSG머니 최대 5% 적립'', ''장바구니 최대 10% 할인쿠폰'', ''스마일클럽 단독 혜택 - 스타벅스 상품 전용 딜'', ''매월 4장씩 최대 12% 할인쿠폰'', ''스마일배송 1만 5천원 이상 무제한 무료배송'', ''스마일배송 상품 스마일캐시 1% 적립'', ''카드 할인 혜택 자세히보기 (레이어팝업 열기)'', ''2022.10.26부터, SSG.COM 삼성카드 발급 시 스마일클럽에 자동 가입됩니다.'', ''G마켓 또는 옥션을 통해 스마일클럽에 가입한 경우 혜택을 받을 수 없습니다.'', ''SSG.COM 삼성카드로 스마일클럽 월 이용료(3,900원) 정기결제 시 혜택이 제공됩니다.(월 1회)'', ''SSG.COM을 통한 스마일클럽 가입건에 한해 혜택이 제공됩니다.'', ''결제금액이 할인금액보다 적을 경우, 결제금액만큼 할인이 적용됩니다.(결제금액이 없는 경우 할인 대상에서 제외)'', ''01. 5만 5천원 이상 결제 시 사용 가능한 5만원 할인쿠폰 제공'', ''50,000원 할인쿠폰 - SSGPAY 바로결제 이용 시 사용 가능'', ''쿠폰 발급기간 : 2023.03.01 ~ 2023.03.31'', ''쿠폰 사용기간 : 2023.03.01 ~ 2023.03.31'', ''SSGPAY 바로결제 이용 이력이 없는 회원'', ''SSGPAY 바로결제에 등록된 SSG.COM 삼성카드로 결제해야 쿠폰을 사용할 수 있습니다.'', ''배송비 등을 제외한 최종 결제금액이 55,000원 이상이여야 쿠폰을 사용할 수 있습니다.'', ''SSGPAY 바로결제에 등록된 SSG.COM 삼성카드로 결제하셔야 쿠폰을 사용할 수 있습니다.'', ''쿠폰은 통합 회원 본인 명의의 SSG.COM 삼성카드로 결제 시 사용 가능하며, 다른 부정적인 방법으로 사용한 경우에는 주문이 취소될 수 있습니다.'', ''쿠폰은 결제 화면에서 자동으로 적용됩니다.'', ''결제 화면에서 쿠폰 변경을 원할 경우 ‘쿠폰선택’을 눌러주세요.'', ''03. 스마일클럽 월 이용료 결제 카드로 SSG.COM 삼성카드 등록 시 SSG머니 3,900원 적립'', ''행사기간 동안 SSG.COM 삼성카드를 통해 스마일클럽 자동 가입 시 SSG머니 3,900원 즉시 적립'', ''적립시점까지 스마일클럽 월 이용료 정기결제 수단에 SSG.COM 삼성카드를 등록해야 혜택을 받을 수 있습니다.'', ''본 상품 발급 및 SSGPAY앱 설치시 바로결제에 자동 등록 됩니다.'', ''SSG MONEY는 매월 1일~말일까지 매출전표가 접수된 금액에 대해 다음달 25일 SSG.COM 계정으로 적립됩니다.'', ''적립된 SSG MONEY는 SSGPAY 회원가입 후 조회 및 사용 가능합니다. 단, SSG.COM 에서는 SSG.COM만 가입해도 사용 가능합니다.'', ''연체이자율 : 회원별/이용상품별 정상이자율+3.0%p(최고 연 20.0%)'', ''이미 SSG.COM 삼성카드를 가지고 계시네요!스마일클럽 가입하고 모든 혜택 누리세요'', ''이미 SSG.COM 삼성카드를 가지고 계시네요!G마켓 또는 옥션을 통해 스마일클럽에 가입한 경우, 월 이용료 결제일할인 혜택을 받을 수 없습니다.'', ''쿠폰 다운되었습니다.다운받은 쿠폰은 SSG.COM 삼성카드 발급 후 바로 사용 가능합니다.'', ''이벤트 쿠폰이 이미 발급되었습니다.다운받은 쿠폰은 SSG.COM 삼성카드 발급 후 바로 사용 가능합니다.'']'

# Row 8
$ws.Range("A8").Value = 'SK-II 3/2(목) 8PM'
$ws.Range("B8").Value = '극광 에센스 리뉴얼 최초 공개! 선착순 핫딜 20만원대+백화점 상품권 증정+선물포장'
$ws.Range("C8").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000004573&domainSiteNo=6005'
$ws.Range("D8").Value = '[SSG.LIVE]3/2(목) SK-II'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '2023-02-23'
$ws.Range("E8").Style = "Normal"
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = '2023-03-02'
$ws.Range("F8").Style = "Normal"
$ws.Range("G8").Value = '[''이벤트/쿠폰 > [SSG.LIVE]3/2(목) SK-II'', ''스마일클럽'', ''- 사은품 지급 및 이벤트 혜택 제공'', ''- 사은품 지급 및 이벤트 혜택 제공 관련 업무 종료 후 즉시 파기'']'

# Row 9
$ws.Range("A9").Value = '보라카이(하나투어) 3/2(목) 9PM'
$ws.Range("B9").Value = '보라카이 베스트셀러 헤난가든 초특가 4/5일'
$ws.Range("C9").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000004579&siteNo=6005'
$ws.Range("D9").Value = '보라카이 여행(하나투어) @SSG.LIVE 3/2(목) 9PM'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '2023-02-21'
$ws.Range("E9").Style = "Normal"
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = '2023-03-02'
$ws.Range("F9").Style = "Normal"
$ws.Range("G9").Value = '[''이벤트/쿠폰 > 보라카이 여행(하나투어) @SSG.LIVE 3/2(목) 9PM'', ''스마일클럽'', ''이벤트 혜택 당첨 주의사항'', ''- 사은품 지급 및 이벤트 혜택 제공'', ''- 사은품 지급 및 이벤트 혜택 제공 관련 업무 종료 후 즉시 파기'']'

# Row 10
$ws.Range("A10").Value = '3월 맘키즈 플러스'
$ws.Range("B10").Value = '지금 할인 중! ~40% 할인 혜택'
$ws.Range("C10").Value = 'https://event.ssg.com/eventDetail.ssg?nevntId=1000000001665'
$ws.Range("D10").Value = '이달의 맘키즈 PLUS'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '2022-09-01'
$ws.Range("E10").Style = "Normal"
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = '2999-12-13'
$ws.Range("F10").Style = "Normal"
$ws.Range("G10").Value = '[''이벤트/쿠폰 > 이달의 맘키즈 PLUS'', ''스마일클럽'', ''맘키즈 ~40% 쿠폰상품'', ''지금 할인 중!'', ''※ 할인 금액은 상품 별로 상이하니 각 상품페이지를 꼭 참조하세요'', ''맘키즈 에누리 쿠폰'', ''행사대상 상품별 최대 40% 할인'', ''맘키즈 상품 에누리 쿠폰'', ''맘키즈 클럽 회원이라면 로그인 후 각 상품 상세페이지에서도 쿠폰을 받으실 수 있습니다.'', ''쿠폰이 적용되지 않는다면?'', ''맘키즈 에누리 쿠폰 ~40%'', ''5無키즈 100% 유기농 망고오렌지 500ml (100mlx5포)'', ''5無 100% 국산 배도라지즙 400ml (100mlx4)'', ''시크릿쥬쥬 베스트프렌즈마고 (N2쓱배송, 전국택배)'', ''시크릿쥬쥬 베스트프렌즈스쿨백세트 (N2쓱배송, 전국택배)'', ''시크릿쥬쥬 베스트프렌즈리나 (N2쓱배송, 전국택배)'', ''시크릿쥬쥬 베스트프렌즈쥬쥬 (N2쓱배송, 전국택배)'', ''베스트프렌즈시크릿프렌즈폰 (N2쓱배송, 전국택배)'', ''시크릿쥬쥬 베스트프렌즈미니사물함 (N2쓱배송, 전국택배)'', ''★금주의 할인딜★ KIDS 침구 / 어린이집 낮잠이불'']'

# Remove the now-unused last row (11) so the used range shrinks to A1:G10
$ws.Rows.Item(11).Delete()
